# "add proper mods and extremeties to dodge"
#
# Sheet2 drives an INDEX/MATCH code-generator off dropdown selections that
# mirror Sheet1's lookup tables. This fixes up the "Attack" row's modifiers
# (Melee Weapon -> All, Crit Threshold -> Roll Add, Critical -> Targetted)
# and adds a proper modifier for the "Defence / Full Defense Dodge" row
# (None -> Roll Add), then leaves the selection sitting on Sheet2 where the
# work happened.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1 was the active sheet/selection before; move its cursor even though
# it's no longer the selected tab.
$ws1.Activate()
$ws1.Range("G8").Select()

# Sheet2 becomes the active tab, with new dropdown picks driving the
# INDEX/MATCH formulas in row 4 (Attack/Weapon) and row 6 (Defence/Full
# Defense Dodge extremity modifier).
$ws2.Activate()

$ws2.Range("C3").Value = "All"
$ws2.Range("E3").Value = "Roll Add"
$ws2.Range("F3").Value = "Targetted"

$ws2.Range("C5").Value = "Roll Add"

$ws2.Range("I4").Select()
